# Update the "Förändrad" (Changed) date column (C) for data rows 2-29
# from 2024-10-06 (serial 45571) to 2024-10-07 (serial 45572).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45571) {
        $cell.Value2 = 45572
    }
}
